$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so numeric-looking strings
# (e.g. "0.619", "246.38") are not reinterpreted as numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '42.104.51'
$ws.Range("E2").Value = '  -2.08%  '

$ws.Range("D3").Value = '2.243.86'
$ws.Range("E3").Value = '  -2.45%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").Value = '246.38'
$ws.Range("E5").Value = '  -2.14%  '

$ws.Range("D6").Value = '0.619'
$ws.Range("E6").Value = '  -4.90%  '

$ws.Range("D7").Value = '73.95'
$ws.Range("E7").Value = '  -1.73%  '

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("D9").Value = '0.613'
$ws.Range("E9").Value = '  -5.75%  '

$ws.Range("D10").Value = '41.37'
$ws.Range("E10").Value = '  +5.41%  '

$ws.Range("D11").Value = '0.0934'
$ws.Range("E11").Value = '  -5.37%  '

$ws.Range("D12").Value = '7.08'
$ws.Range("E12").Value = '  -6.32%  '

$ws.Range("E13").Value = '  -4.44%  '

$ws.Range("D14").Value = '2.578.51'
$ws.Range("E14").Value = '  -2.50%  '

$ws.Range("D15").Value = '14.41'
$ws.Range("E15").Value = '  -5.28%  '

$ws.Range("D16").Value = '0.847'
$ws.Range("E16").Value = '  -3.84%  '

$ws.Range("D17").Value = '2.236.20'
$ws.Range("E17").Value = '  -2.81%  '

$ws.Range("D18").Value = '42.084.51'
$ws.Range("E18").Value = '  -1.89%  '

$ws.Range("D19").Value = '0.0₃0969'
$ws.Range("E19").Value = '  -3.57%  '

$ws.Range("D20").Value = '71.64'
$ws.Range("E20").Value = '  -1.10%  '

$ws.Range("D21").Value = '6.07'
$ws.Range("E21").Value = '  -3.21%  '

$ws.Range("D22").Value = '2.34'
$ws.Range("E22").Value = '  +6.73%  '

$ws.Range("D23").Value = '228.93'
$ws.Range("E23").Value = '  -3.61%  '

$ws.Range("E24").Value = '  +0.07%  '

$ws.Range("D25").Value = '10.96'
$ws.Range("E25").Value = '  -4.08%  '

$ws.Range("D26").Value = '3.54'
$ws.Range("E26").Value = '  -8.76%  '

$ws.Range("D27").Value = '2.29'
$ws.Range("E27").Value = '  -4.04%  '

$ws.Range("D28").Value = '7.37'
$ws.Range("E28").Value = '  +17.57%  '

$ws.Range("D29").Value = '2.24'
$ws.Range("E29").Value = '  +2.74%  '

$ws.Range("D30").Value = '169.62'
$ws.Range("E30").Value = '  +0.99%  '

$ws.Range("D31").Value = '20.62'
$ws.Range("E31").Value = '  -2.17%  '

$ws.Range("D32").Value = '0.0830'
$ws.Range("E32").Value = '  -3.90%  '

$ws.Range("D33").Value = '0.120'
$ws.Range("E33").Value = '  -6.97%  '

$ws.Range("D34").Value = '30.18'
$ws.Range("E34").Value = '  -3.41%  '

$ws.Range("D35").Value = '0.125'
$ws.Range("E35").Value = '  -3.07%  '

$ws.Range("D36").Value = '4.48'
$ws.Range("E36").Value = '  -4.29%  '

$ws.Range("D37").Value = '4.82'
$ws.Range("E37").Value = '  +0.20%  '

$ws.Range("D38").Value = '0.0298'
$ws.Range("E38").Value = '  -2.49%  '

$ws.Range("D39").Value = '13.26'
$ws.Range("E39").Value = '  -3.57%  '

$ws.Range("D40").Value = '2.17'
$ws.Range("E40").Value = '  -6.26%  '

$ws.Range("D41").Value = '5.77'
$ws.Range("E41").Value = '  -4.17%  '

$ws.Range("D42").Value = '110.04'
$ws.Range("E42").Value = '  +4.13%  '

$ws.Range("D43").Value = '0.202'
$ws.Range("E43").Value = '  -4.50%  '

$ws.Range("D44").Value = '61.06'
$ws.Range("E44").Value = '  -0.48%  '

$ws.Range("D45").Value = '8.67'
$ws.Range("E45").Value = '  -4.91%  '

$ws.Range("D46").Value = '0.1000'
$ws.Range("E46").Value = '  -1.39%  '

$ws.Range("D47").Value = '0.997'
$ws.Range("E47").Value = '  -0.33%  '

$ws.Range("D48").Value = '1.11'
$ws.Range("E48").Value = '  -5.13%  '

$ws.Range("D49").Value = '1.16'
$ws.Range("E49").Value = '  -2.17%  '

$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").Value = '2.25'
$ws.Range("E50").Value = '  -1.94%  '

$ws.Range("E51").Value = '  -1.22%  '

# Drop the temporary text-number-format so no stray style lingers
# on the price column (matches original unstyled cells).
$ws.Range("D2:D51").ClearFormats()
